# Apply the "Updated cryptos list" edit: refresh Price (D) and Volume(1h) (E) columns.
#
# D-column values are typed as text in the source file (e.g. "56.532.45", "0.999",
# "1.00"), but Range.Value / Range.Formula auto-convert plain decimal-looking
# strings into real numbers (losing the exact text, e.g. "1.00" -> 1). To keep the
# cell a genuine text value we enter it as a text-literal formula (="123") and then
# Copy / PasteSpecial(values) it onto itself, which flattens the formula to its
# computed (text) result without touching number formatting/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = '="56.532.45"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E2").Value = '  +3.68%  '

$c = $ws.Range("D3")
$c.Formula = '="2.994.83"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E3").Value = '  +4.06%  '

$c = $ws.Range("D4")
$c.Formula = '="0.999"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E4").Value = '  +0.05%  '

$c = $ws.Range("D5")
$c.Formula = '="507.23"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E5").Value = '  +8.14%  '

$c = $ws.Range("D6")
$c.Formula = '="138.26"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E6").Value = '  +10.14%  '

$c = $ws.Range("D8")
$c.Formula = '="0.432"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E8").Value = '  +6.45%  '

$c = $ws.Range("D9")
$c.Formula = '="7.56"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E9").Value = '  +14.40%  '

$ws.Range("E10").Value = '  +11.40%  '

$c = $ws.Range("D11")
$c.Formula = '="0.351"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E11").Value = '  +5.92%  '

$ws.Range("E12").Value = '  +5.44%  '

$c = $ws.Range("D13")
$c.Formula = '="3.499.31"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E13").Value = '  +3.90%  '

$c = $ws.Range("D14")
$c.Formula = '="25.39"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E14").Value = '  +8.67%  '

$c = $ws.Range("D15")
$c.Formula = '="0.0000153"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E15").Value = '  +14.70%  '

$c = $ws.Range("D16")
$c.Formula = '="56.547.83"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E16").Value = '  +3.94%  '

$c = $ws.Range("D17")
$c.Formula = '="2.992.77"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E17").Value = '  +4.27%  '

$c = $ws.Range("D18")
$c.Formula = '="5.86"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E18").Value = '  +9.67%  '

$c = $ws.Range("D19")
$c.Formula = '="12.40"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E19").Value = '  +8.44%  '

$c = $ws.Range("D20")
$c.Formula = '="7.82"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E20").Value = '  +10.74%  '

$c = $ws.Range("D21")
$c.Formula = '="326.81"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E21").Value = '  +8.80%  '

$c = $ws.Range("D22")
$c.Formula = '="1.00"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E22").Value = '  +0.05%  '

$c = $ws.Range("D23")
$c.Formula = '="0.479"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E23").Value = '  +9.12%  '

$c = $ws.Range("D24")
$c.Formula = '="62.37"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E24").Value = '  +6.00%  '

$c = $ws.Range("D25")
$c.Formula = '="0.171"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E25").Value = '  +14.18%  '

$c = $ws.Range("D26")
$c.Formula = '="0.999"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E26").Value = '  +0.09%  '

$c = $ws.Range("D27")
$c.Formula = '="0.0₃0905"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E27").Value = '  +13.44%  '

$c = $ws.Range("D28")
$c.Formula = '="6.60"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E28").Value = '  +7.63%  '

$c = $ws.Range("D29")
$c.Formula = '="7.07"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E29").Value = '  +13.98%  '

$c = $ws.Range("D30")
$c.Formula = '="1.27"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E30").Value = '  +15.34%  '

$c = $ws.Range("D31")
$c.Formula = '="1.78"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E31").Value = '  +11.11%  '

$c = $ws.Range("D32")
$c.Formula = '="20.59"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E32").Value = '  +9.41%  '

$c = $ws.Range("D33")
$c.Formula = '="155.82"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E33").Value = '  +10.06%  '

$c = $ws.Range("D34")
$c.Formula = '="4.50"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E34").Value = '  +8.09%  '

$c = $ws.Range("D35")
$c.Formula = '="5.62"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E35").Value = '  +4.21%  '

$c = $ws.Range("D36")
$c.Formula = '="1.26"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E36").Value = '  +4.31%  '

$c = $ws.Range("D37")
$c.Formula = '="0.0679"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E37").Value = '  +9.71%  '

$c = $ws.Range("D38")
$c.Formula = '="23.87"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E38").Value = '  +4.55%  '

$c = $ws.Range("D39")
$c.Formula = '="3.024.67"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E39").Value = '  +4.07%  '

$c = $ws.Range("D40")
$c.Formula = '="36.91"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E40").Value = '  +5.17%  '

$c = $ws.Range("D41")
$c.Formula = '="0.999"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E41").Value = '  +0.06%  '

$c = $ws.Range("D42")
$c.Formula = '="0.645"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E42").Value = '  +6.92%  '

$c = $ws.Range("D43")
$c.Formula = '="2.258.10"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E43").Value = '  +11.23%  '

$c = $ws.Range("D44")
$c.Formula = '="1.41"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E44").Value = '  +7.71%  '

$c = $ws.Range("D45")
$c.Formula = '="0.988"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E45").Value = '  +4.87%  '

$c = $ws.Range("D46")
$c.Formula = '="3.61"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E46").Value = '  +5.40%  '

$c = $ws.Range("D47")
$c.Formula = '="1.99"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E47").Value = '  +25.77%  '

$c = $ws.Range("D48")
$c.Formula = '="0.0237"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E48").Value = '  +11.08%  '

$c = $ws.Range("D49")
$c.Formula = '="5.76"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E49").Value = '  +8.00%  '

$c = $ws.Range("D50")
$c.Formula = '="19.11"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E50").Value = '  +8.12%  '

$ws.Range("E51").Value = '  +9.53%  '

$excel.CutCopyMode = 0
